$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F column) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 7213
$wsExhibit.Range("F12").Value = 218
$wsExhibit.Range("F16").Value = 1856
$wsExhibit.Range("F18").Value = 36

# Sheet "全部类型" (sheet4): update 想去人数 (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7213
$wsAll.Range("F13").Value = 218
$wsAll.Range("F17").Value = 1856
$wsAll.Range("F19").Value = 36
